# Add a new "2022" column (K) to the table, mirroring the formatting of
# the existing "2021" column (J), then move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column J's formatting (styles) for rows 4-14 onto column K so the
# new cells pick up the same number formats / fonts / borders as the rest
# of the table (xlPasteFormats = -4122).
$ws.Range("J4:J14").Copy() | Out-Null
$ws.Range("K4:K14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# New "2022" year header + data values.
$ws.Range("K4").Value = 2022
$ws.Range("K5").Value = 26.495524312074597
$ws.Range("K6").Value = 59.383769502755833
$ws.Range("K7").Value = 38.32334404557426
$ws.Range("K8").Value = 48.136790950525594
$ws.Range("K9").Value = 46.63213064070051
$ws.Range("K10").Value = 32.657429481680126
$ws.Range("K11").Value = 31.457245964894081
$ws.Range("K12").Value = 22.734405597714229
$ws.Range("K13").Value = -0.19691879995369213
$ws.Range("K14").Value = 33.158040409631916

# Update the saved selection/active cell shown in the sheet view.
$ws.Range("M7").Select() | Out-Null
